$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 41.428665
$ws.Range("H2").Value2 = 124.285995
$ws.Range("I2").Value2 = 0.06969137269740189
$ws.Range("J2").Value2 = 0.06969137269740189
$ws.Range("M2").Value2 = 6.305846
$ws.Range("N2").Value2 = 18.917538
$ws.Range("O2").Value2 = 0.01356150511917599
$ws.Range("P2").Value2 = 0.01356150511917599
$ws.Range("Q2").Value2 = 261.2427814755899
$ws.Range("R2").Value2 = 2351.18503328031
$ws.Range("S2").Value2 = 0.0009451199075982174
$ws.Range("T2").Value2 = 0.0009451199075982173
$ws.Range("G3").Value2 = 41.428665
$ws.Range("H3").Value2 = 124.285995
$ws.Range("I3").Value2 = 0.06969137269740189
$ws.Range("J3").Value2 = 0.06969137269740189
$ws.Range("O3").Value2 = 0.392557056479861
$ws.Range("P3").Value2 = 0.3925570564798609
$ws.Range("Q3").Value2 = 7562.043919274085
$ws.Range("R3").Value2 = 68058.39527346675
$ws.Range("S3").Value2 = 0.02735784012813304
$ws.Range("T3").Value2 = 0.02735784012813303
$ws.Range("G4").Value2 = 41.428665
$ws.Range("H4").Value2 = 124.285995
$ws.Range("I4").Value2 = 0.06969137269740189
$ws.Range("J4").Value2 = 0.06969137269740189
$ws.Range("M4").Value2 = 127.396393
$ws.Range("N4").Value2 = 382.189179
$ws.Range("O4").Value2 = 0.2739817680029065
$ws.Range("P4").Value2 = 0.2739817680029065
$ws.Range("Q4").Value2 = 5277.862487805344
$ws.Range("R4").Value2 = 47500.7623902481
$ws.Range("S4").Value2 = 0.01909416550618366
$ws.Range("T4").Value2 = 0.01909416550618366
$ws.Range("G5").Value2 = 41.428665
$ws.Range("H5").Value2 = 124.285995
$ws.Range("I5").Value2 = 0.06969137269740189
$ws.Range("J5").Value2 = 0.06969137269740189
$ws.Range("M5").Value2 = 19.42400133333333
$ws.Range("N5").Value2 = 58.272004
$ws.Range("O5").Value2 = 0.04177372766745037
$ws.Range("P5").Value2 = 0.04177372766745036
$ws.Range("Q5").Value2 = 804.7104441982199
$ws.Range("R5").Value2 = 7242.393997783979
$ws.Range("S5").Value2 = 0.002911268423832053
$ws.Range("T5").Value2 = 0.002911268423832052
$ws.Range("G6").Value2 = 41.428665
$ws.Range("H6").Value2 = 124.285995
$ws.Range("I6").Value2 = 0.06969137269740189
$ws.Range("J6").Value2 = 0.06969137269740189
$ws.Range("M6").Value2 = 129.3233566666667
$ws.Range("N6").Value2 = 387.97007
$ws.Range("O6").Value2 = 0.2781259427306063
$ws.Range("P6").Value2 = 0.2781259427306062
$ws.Range("Q6").Value2 = 5357.69402001885
$ws.Range("R6").Value2 = 48219.24618016965
$ws.Range("S6").Value2 = 0.01938297873165493
$ws.Range("T6").Value2 = 0.01938297873165493
$ws.Range("H7").Value2 = 510.696747
$ws.Range("I7").Value2 = 0.2863649869040173
$ws.Range("J7").Value2 = 0.2863649869040173
$ws.Range("M7").Value2 = 6.305846
$ws.Range("N7").Value2 = 18.917538
$ws.Range("O7").Value2 = 0.01356150511917599
$ws.Range("P7").Value2 = 0.01356150511917599
$ws.Range("Q7").Value2 = 1073.458346427654
$ws.Range("R7").Value2 = 9661.125117848886
$ws.Range("S7").Value2 = 0.003883540235851595
$ws.Range("T7").Value2 = 0.003883540235851596
$ws.Range("H8").Value2 = 510.696747
$ws.Range("I8").Value2 = 0.2863649869040173
$ws.Range("J8").Value2 = 0.2863649869040173
$ws.Range("O8").Value2 = 0.392557056479861
$ws.Range("P8").Value2 = 0.3925570564798609
$ws.Range("Q8").Value2 = 31072.7787973569
$ws.Range("R8").Value2 = 279655.0091762121
$ws.Range("S8").Value2 = 0.112414596337935
$ws.Range("T8").Value2 = 0.112414596337935
$ws.Range("H9").Value2 = 510.696747
$ws.Range("I9").Value2 = 0.2863649869040173
$ws.Range("J9").Value2 = 0.2863649869040173
$ws.Range("M9").Value2 = 127.396393
$ws.Range("N9").Value2 = 382.189179
$ws.Range("O9").Value2 = 0.2739817680029065
$ws.Range("P9").Value2 = 0.2739817680029065
$ws.Range("Q9").Value2 = 21686.97449487786
$ws.Range("R9").Value2 = 195182.7704539007
$ws.Range("S9").Value2 = 0.07845878540609183
$ws.Range("T9").Value2 = 0.07845878540609184
$ws.Range("H10").Value2 = 510.696747
$ws.Range("I10").Value2 = 0.2863649869040173
$ws.Range("J10").Value2 = 0.2863649869040173
$ws.Range("M10").Value2 = 19.42400133333333
$ws.Range("N10").Value2 = 58.272004
$ws.Range("O10").Value2 = 0.04177372766745037
$ws.Range("P10").Value2 = 0.04177372766745036
$ws.Range("Q10").Value2 = 3306.591431552332
$ws.Range("R10").Value2 = 29759.32288397099
$ws.Range("S10").Value2 = 0.01196253297642141
$ws.Range("T10").Value2 = 0.01196253297642141
$ws.Range("H11").Value2 = 510.696747
$ws.Range("I11").Value2 = 0.2863649869040173
$ws.Range("J11").Value2 = 0.2863649869040173
$ws.Range("M11").Value2 = 129.3233566666667
$ws.Range("N11").Value2 = 387.97007
$ws.Range("O11").Value2 = 0.2781259427306063
$ws.Range("P11").Value2 = 0.2781259427306062
$ws.Range("Q11").Value2 = 22015.00585359581
$ws.Range("R11").Value2 = 198135.0526823623
$ws.Range("S11").Value2 = 0.07964553194771752
$ws.Range("T11").Value2 = 0.07964553194771752
$ws.Range("G12").Value2 = 244.5761666666666
$ws.Range("H12").Value2 = 733.7284999999999
$ws.Range("I12").Value2 = 0.4114264551867299
$ws.Range("J12").Value2 = 0.41142645518673
$ws.Range("M12").Value2 = 6.305846
$ws.Range("N12").Value2 = 18.917538
$ws.Range("O12").Value2 = 0.01356150511917599
$ws.Range("P12").Value2 = 0.01356150511917599
$ws.Range("Q12").Value2 = 1542.259642270333
$ws.Range("R12").Value2 = 13880.336780433
$ws.Range("S12").Value2 = 0.005579561978179269
$ws.Range("T12").Value2 = 0.005579561978179269
$ws.Range("G13").Value2 = 244.5761666666666
$ws.Range("H13").Value2 = 733.7284999999999
$ws.Range("I13").Value2 = 0.4114264551867299
$ws.Range("J13").Value2 = 0.41142645518673
$ws.Range("O13").Value2 = 0.392557056479861
$ws.Range("P13").Value2 = 0.3925570564798609
$ws.Range("Q13").Value2 = 44642.89916030439
$ws.Range("R13").Value2 = 401786.0924427395
$ws.Range("S13").Value2 = 0.1615083582060461
$ws.Range("T13").Value2 = 0.1615083582060461
$ws.Range("G14").Value2 = 244.5761666666666
$ws.Range("H14").Value2 = 733.7284999999999
$ws.Range("I14").Value2 = 0.4114264551867299
$ws.Range("J14").Value2 = 0.41142645518673
$ws.Range("M14").Value2 = 127.396393
$ws.Range("N14").Value2 = 382.189179
$ws.Range("O14").Value2 = 0.2739817680029065
$ws.Range("P14").Value2 = 0.2739817680029065
$ws.Range("Q14").Value2 = 31158.12144710016
$ws.Range("R14").Value2 = 280423.0930239015
$ws.Range("S14").Value2 = 0.1127233475952289
$ws.Range("T14").Value2 = 0.1127233475952289
$ws.Range("G15").Value2 = 244.5761666666666
$ws.Range("H15").Value2 = 733.7284999999999
$ws.Range("I15").Value2 = 0.4114264551867299
$ws.Range("J15").Value2 = 0.41142645518673
$ws.Range("M15").Value2 = 19.42400133333333
$ws.Range("N15").Value2 = 58.272004
$ws.Range("O15").Value2 = 0.04177372766745037
$ws.Range("P15").Value2 = 0.04177372766745036
$ws.Range("Q15").Value2 = 4750.647787434888
$ws.Range("R15").Value2 = 42755.830086914
$ws.Range("S15").Value2 = 0.01718681669415493
$ws.Range("T15").Value2 = 0.01718681669415493
$ws.Range("G16").Value2 = 244.5761666666666
$ws.Range("H16").Value2 = 733.7284999999999
$ws.Range("I16").Value2 = 0.4114264551867299
$ws.Range("J16").Value2 = 0.41142645518673
$ws.Range("M16").Value2 = 129.3233566666667
$ws.Range("N16").Value2 = 387.97007
$ws.Range("O16").Value2 = 0.2781259427306063
$ws.Range("P16").Value2 = 0.2781259427306062
$ws.Range("Q16").Value2 = 31629.41083399945
$ws.Range("R16").Value2 = 284664.697505995
$ws.Range("S16").Value2 = 0.1144283707131208
$ws.Range("T16").Value2 = 0.1144283707131208
$ws.Range("G17").Value2 = 24.173247
$ws.Range("H17").Value2 = 72.51974100000001
$ws.Range("I17").Value2 = 0.04066427836821081
$ws.Range("J17").Value2 = 0.04066427836821081
$ws.Range("M17").Value2 = 6.305846
$ws.Range("N17").Value2 = 18.917538
$ws.Range("O17").Value2 = 0.01356150511917599
$ws.Range("P17").Value2 = 0.01356150511917599
$ws.Range("Q17").Value2 = 152.432772901962
$ws.Range("R17").Value2 = 1371.894956117658
$ws.Range("S17").Value2 = 0.0005514688192580884
$ws.Range("T17").Value2 = 0.0005514688192580883
$ws.Range("G18").Value2 = 24.173247
$ws.Range("H18").Value2 = 72.51974100000001
$ws.Range("I18").Value2 = 0.04066427836821081
$ws.Range("J18").Value2 = 0.04066427836821081
$ws.Range("O18").Value2 = 0.392557056479861
$ws.Range("P18").Value2 = 0.3925570564798609
$ws.Range("Q18").Value2 = 4412.383442369204
$ws.Range("R18").Value2 = 39711.45098132284
$ws.Range("S18").Value2 = 0.01596304942010252
$ws.Range("T18").Value2 = 0.01596304942010252
$ws.Range("G19").Value2 = 24.173247
$ws.Range("H19").Value2 = 72.51974100000001
$ws.Range("I19").Value2 = 0.04066427836821081
$ws.Range("J19").Value2 = 0.04066427836821081
$ws.Range("M19").Value2 = 127.396393
$ws.Range("N19").Value2 = 382.189179
$ws.Range("O19").Value2 = 0.2739817680029065
$ws.Range("P19").Value2 = 0.2739817680029065
$ws.Range("Q19").Value2 = 3079.584474898072
$ws.Range("R19").Value2 = 27716.26027408264
$ws.Range("S19").Value2 = 0.01114127088188475
$ws.Range("T19").Value2 = 0.01114127088188475
$ws.Range("G20").Value2 = 24.173247
$ws.Range("H20").Value2 = 72.51974100000001
$ws.Range("I20").Value2 = 0.04066427836821081
$ws.Range("J20").Value2 = 0.04066427836821081
$ws.Range("M20").Value2 = 19.42400133333333
$ws.Range("N20").Value2 = 58.272004
$ws.Range("O20").Value2 = 0.04177372766745037
$ws.Range("P20").Value2 = 0.04177372766745036
$ws.Range("Q20").Value2 = 469.541181958996
$ws.Range("R20").Value2 = 4225.870637630965
$ws.Range("S20").Value2 = 0.001698698490347032
$ws.Range("T20").Value2 = 0.001698698490347031
$ws.Range("G21").Value2 = 24.173247
$ws.Range("H21").Value2 = 72.51974100000001
$ws.Range("I21").Value2 = 0.04066427836821081
$ws.Range("J21").Value2 = 0.04066427836821081
$ws.Range("M21").Value2 = 129.3233566666667
$ws.Range("N21").Value2 = 387.97007
$ws.Range("O21").Value2 = 0.2781259427306063
$ws.Range("P21").Value2 = 0.2781259427306062
$ws.Range("Q21").Value2 = 3126.165443572431
$ws.Range("R21").Value2 = 28135.48899215188
$ws.Range("S21").Value2 = 0.01130979075661843
$ws.Range("T21").Value2 = 0.01130979075661843
$ws.Range("G22").Value2 = 114.0486906666667
$ws.Range("H22").Value2 = 342.146072
$ws.Range("I22").Value2 = 0.19185290684364
$ws.Range("J22").Value2 = 0.19185290684364
$ws.Range("M22").Value2 = 6.305846
$ws.Range("N22").Value2 = 18.917538
$ws.Range("O22").Value2 = 0.01356150511917599
$ws.Range("P22").Value2 = 0.01356150511917599
$ws.Range("Q22").Value2 = 719.1734798456373
$ws.Range("R22").Value2 = 6472.561318610737
$ws.Range("S22").Value2 = 0.002601814178288818
$ws.Range("T22").Value2 = 0.002601814178288817
$ws.Range("G23").Value2 = 114.0486906666667
$ws.Range("H23").Value2 = 342.146072
$ws.Range("I23").Value2 = 0.19185290684364
$ws.Range("J23").Value2 = 0.19185290684364
$ws.Range("O23").Value2 = 0.392557056479861
$ws.Range("P23").Value2 = 0.3925570564798609
$ws.Range("Q23").Value2 = 20817.4993753006
$ws.Range("R23").Value2 = 187357.4943777054
$ws.Range("S23").Value2 = 0.0753132123876443
$ws.Range("T23").Value2 = 0.07531321238764428
$ws.Range("G24").Value2 = 114.0486906666667
$ws.Range("H24").Value2 = 342.146072
$ws.Range("I24").Value2 = 0.19185290684364
$ws.Range("J24").Value2 = 0.19185290684364
$ws.Range("M24").Value2 = 127.396393
$ws.Range("N24").Value2 = 382.189179
$ws.Range("O24").Value2 = 0.2739817680029065
$ws.Range("P24").Value2 = 0.2739817680029065
$ws.Range("Q24").Value2 = 14529.3918173061
$ws.Range("R24").Value2 = 130764.5263557549
$ws.Range("S24").Value2 = 0.0525641986135174
$ws.Range("T24").Value2 = 0.0525641986135174
$ws.Range("G25").Value2 = 114.0486906666667
$ws.Range("H25").Value2 = 342.146072
$ws.Range("I25").Value2 = 0.19185290684364
$ws.Range("J25").Value2 = 0.19185290684364
$ws.Range("M25").Value2 = 19.42400133333333
$ws.Range("N25").Value2 = 58.272004
$ws.Range("O25").Value2 = 0.04177372766745037
$ws.Range("P25").Value2 = 0.04177372766745036
$ws.Range("Q25").Value2 = 2215.281919574255
$ws.Range("R25").Value2 = 19937.53727616829
$ws.Range("S25").Value2 = 0.008014411082694942
$ws.Range("T25").Value2 = 0.00801441108269494
$ws.Range("G26").Value2 = 114.0486906666667
$ws.Range("H26").Value2 = 342.146072
$ws.Range("I26").Value2 = 0.19185290684364
$ws.Range("J26").Value2 = 0.19185290684364
$ws.Range("M26").Value2 = 129.3233566666667
$ws.Range("N26").Value2 = 387.97007
$ws.Range("O26").Value2 = 0.2781259427306063
$ws.Range("P26").Value2 = 0.2781259427306062
$ws.Range("Q26").Value2 = 14749.15950045167
$ws.Range("R26").Value2 = 132742.435504065
$ws.Range("S26").Value2 = 0.05335927058149455
$ws.Range("T26").Value2 = 0.05335927058149454
